$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Bugs and errors"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Bugs and errors")

# Row 9: the Git-branch marker in column I is cleared (style is kept)
$ws1.Cells.Item(9, 9).Value = ""

# Row 11: new bug entry
#  - fill "found by" / "resolved by" (Ton) before the description so the
#    shared-string table order matches the source edit
$ws1.Cells.Item(11, 3).Value = "Ton"
$ws1.Cells.Item(11, 2).Value = "f.calc_omega uses 3x1 vectors"
$ws1.Cells.Item(11, 5).Value = "Ton"
$ws1.Cells.Item(11, 6).Value = "fixed vector notation and associated matrix math"

$ws1.Cells.Item(11, 4).Value = 44649
$ws1.Cells.Item(11, 4).NumberFormat = "d-mmm"

$ws1.Cells.Item(11, 7).Value = 44649
$ws1.Cells.Item(11, 7).NumberFormat = "d-mmm"

$ws1.Cells.Item(11, 9).Value = "Thomas_Workspace"

# Restore the selection left on this sheet
$ws1.Activate()
$ws1.Range("B12").Select()

# ---------------------------------------------------------------------------
# Sheet "Implemented Features"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Implemented Features")

# Widen columns D and E (best effort - ColumnWidth snaps to a 1/6 character
# grid in Excel, so we pick the closest achievable value)
$ws2.Columns.Item(4).ColumnWidth = 18
$ws2.Columns.Item(5).ColumnWidth = 17.8333333333333

# Existing rows 2-6 gain a "Commited on" date + git branch note
$ws2.Cells.Item(2, 5).Value = 44645
$ws2.Cells.Item(2, 5).NumberFormat = "d-mmm"
$ws2.Cells.Item(2, 6).Value = "Thomas_workspace"

$ws2.Cells.Item(3, 5).Value = 44645
$ws2.Cells.Item(3, 5).NumberFormat = "d-mmm"
$ws2.Cells.Item(3, 6).Value = "Thomas_workspace"

$ws2.Cells.Item(4, 5).Value = 44645
$ws2.Cells.Item(4, 5).NumberFormat = "d-mmm"
$ws2.Cells.Item(4, 6).Value = "Thomas_workspace"

$ws2.Cells.Item(5, 5).Value = 44645
$ws2.Cells.Item(5, 5).NumberFormat = "d-mmm"
$ws2.Cells.Item(5, 6).Value = "Thomas_workspace"

$ws2.Cells.Item(6, 5).Value = 44645
$ws2.Cells.Item(6, 5).NumberFormat = "d-mmm"
$ws2.Cells.Item(6, 6).Value = "Thomas_workspace"

# New row 8: Faulty pitch remover in main_template
$ws2.Cells.Item(8, 1).Value = "Faulty pitch remover in main_template"
$ws2.Cells.Item(8, 2).Value = 44648
$ws2.Cells.Item(8, 3).Value = "Thomas"
$ws2.Cells.Item(8, 5).Value = 44648
$ws2.Cells.Item(8, 5).NumberFormat = "d-mmm"
$ws2.Cells.Item(8, 6).Value = "Thomas_workspace"

# New row 9: Rotation of new optitrack dataset
$ws2.Cells.Item(9, 1).Value = "Rotation of new optitrack dataset"
$ws2.Cells.Item(9, 2).Value = 44643
$ws2.Cells.Item(9, 3).Value = "Thomas"
$ws2.Cells.Item(9, 5).Value = 44643
$ws2.Cells.Item(9, 5).NumberFormat = "d-mmm"
$ws2.Cells.Item(9, 6).Value = "Thomas_workspace"

# Restore the selection left on this sheet
$ws2.Activate()
$ws2.Range("C13").Select()
